$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The post "「我が魂と共に」مع نفسي" (row 774) was removed from the sheet.
# Deleting the entire row shifts every subsequent row up by one, matching
# the diff (old row 775 -> new row 774, ..., old row 797 -> new row 796).
$ws.Rows.Item(774).Delete()
